# "update: reservoir bounds changed to circle with const pressure"
#
# The constraint-settings table switches the reservoir geometry from a
# rectangle (res_width x res_length) to a circle (res_radius), and the
# reservoir outer-boundary condition becomes a fixed constant-pressure
# boundary, so the per-field "boundary_code" choice column is no longer
# needed and is dropped.
#
# Before: field_name | boundary_code | permeability | skin | res_width | res_length | pressure_initial | length_hor_well_bore | length_half_fracture | number_fractures
# After:  field_name | permeability | skin | res_radius | pressure_initial | length_hor_well_bore | length_half_fracture | number_fractures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the "boundary_code" column (column B) entirely - the reservoir
#    boundary condition is now always constant-pressure, so it is no
#    longer a per-row configurable setting.
$ws.Columns.Item(2).Delete()

# After the delete above, the remaining columns shift left by one:
#   A field_name | B permeability | C skin | D res_width | E res_length | F pressure_initial | G length_hor_well_bore | H length_half_fracture | I number_fractures

# 2) Drop the "res_length" column (now column E) - a circular reservoir
#    only needs a single radius, so the separate width/length pair
#    collapses into one dimension.
$ws.Columns.Item(5).Delete()

# After this second delete, columns shift left again:
#   A field_name | B permeability | C skin | D res_width | E pressure_initial | F length_hor_well_bore | G length_half_fracture | H number_fractures

# 3) Repurpose the remaining "res_width" column (D) as "res_radius",
#    keeping its bounds/val_test_period JSON (same [100, 1000] / 500
#    values the width & length columns already shared).
$ws.Range("D1").Value = "res_radius"
$ws.Range("D2").Value = "{""is_discrete"": false, ""bounds"": [100, 1000], ""val_test_period"": 500}"

# Leave the selection on the new res_radius default-bounds cell.
$ws.Range("D2").Select()

Write-Output "Reservoir bounds updated to circle (res_radius) with constant-pressure boundary."
